# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# Only column G ("K") values change; all other columns/rows stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row number -> new K value (column G)
$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    6  = 1
    7  = 2
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 2
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 3
    18 = 1
    19 = 2
    20 = 3
    21 = 4
    22 = 3
    23 = 3
    24 = 3
    25 = 1
    26 = 2
    27 = 1
    28 = 2
    29 = 2
    30 = 4
    31 = 0
    32 = 1
    33 = 2
    34 = 1
    35 = 1
    36 = 1
    38 = 3
    39 = 1
    40 = 2
    41 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
